$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Correcting data analysis": the %-survival formulas in C20:C28 / H20:H28
# used to divide by the cell directly above (B19/G19, then B20/G20, ...).
# That chain divisor is replaced with the fixed Day-0 control reading
# (rounded to 6 d.p., matching the value that was in B19 / G19) so every
# row is normalised against the same baseline instead of compounding.
$ws.Range("C20").Formula = "=(B20/1.125333*C19)"
$ws.Range("H20").Formula = "=G20/1.093667*H19"

$ws.Range("C21").Formula = "=(B21/1.125333*C20)"
$ws.Range("H21").Formula = "=G21/1.093667*H20"

$ws.Range("C22").Formula = "=(B22/1.125333*C21)"
$ws.Range("H22").Formula = "=G22/1.093667*H21"

$ws.Range("C23").Formula = "=(B23/1.125333*C22)"
$ws.Range("H23").Formula = "=G23/1.093667*H22"

$ws.Range("C24").Formula = "=(B24/1.125333*C23)"
$ws.Range("H24").Formula = "=G24/1.093667*H23"

$ws.Range("C25").Formula = "=(B25/1.125333*C24)"
$ws.Range("H25").Formula = "=G25/1.093667*H24"

$ws.Range("C26").Formula = "=(B26/1.125333*C25)"
$ws.Range("H26").Formula = "=G26/1.093667*H25"

$ws.Range("C27").Formula = "=(B27/1.125333*C26)"
$ws.Range("H27").Formula = "=G27/1.093667*H26"

$ws.Range("C28").Formula = "=(B28/1.125333*C27)"
$ws.Range("H28").Formula = "=G28/1.093667*H27"

# "Daily entry": leave the sheet selection where the author finished working
$ws.Range("H20:H28").Select()
